# Weekly update: a new price record for "Comercializadora del Agro de
# Limarí - Arándano (blue)" is inserted at row 6, pushing the previous
# rows 6-12 down to 7-13 (dates/volumes/prices unchanged, just relocated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 6, shifting rows 6:12 down to 7:13.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with this week's record.
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44874
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100101
$ws.Range("H6").Value = "Berries"
$ws.Range("I6").Value = 100101001
$ws.Range("J6").Value = "Arándano (blue)"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 7500
$ws.Range("O6").Value = 8000
$ws.Range("P6").Value = 7750
$ws.Range("Q6").Value = "$/bandeja 2 kilos"
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 3875
$ws.Range("T6").Value = 2
